# Aula 1 - SO
# Sort the roster in A1:A20 alphabetically (A -> Z) and move the
# active selection from A21 to A4, matching the author's re-sort commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sortRange = $ws.Range("A1:A20")

$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1"))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 0
$ws.Sort.Orientation = 1
$ws.Sort.Apply()

$ws.Range("A4").Select()
